$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the two rows that were dropped from the middle of the table
# (old row 71 "-N3ub2xMDKwS4Skdp6kB" and old row 73 "-N3ubhTqk34i3VI-9wIA").
# Deleting row 71 first shifts former row 73 into row 72, so deleting row 72
# next removes it too; all following rows shift up to close the gap.
$ws.Rows.Item(71).Delete()
$ws.Rows.Item(72).Delete()

# Step 2: append the new rows of human-model results to the bottom of the table.
$newRows = @(
  ,@("-N60-G76Th6s4ECC8STb", 22, 20, "closest")
  ,@("-N60-_gZ6U9nuizjY7bg", 7, 22, "random")
  ,@("-N60-i15gPVYo-1ZsNm8", 31, 27, "closest")
  ,@("-N60-v_uOv2Z39ybuhSn", 26, 35, "follow_stag")
  ,@("-N600-7WYheXOKVgUg4K", 18, 29, "closest")
  ,@("-N6003NVcwRQOstmQZia", 26, 20, "closest")
  ,@("-N600GBERJ9DLzOIQDGI", 29, 18, "closest")
  ,@("-N600QSf7tAXOMz9sN03", 9, 20, "random")
  ,@("-N600VFUe-9i_EX2uZc1", 22, 17, "closest")
  ,@("-N600akCmcpWvfUvyXSr", 6, 24, "random")
  ,@("-N600eCV-aGYL8Z_aMNO", 27, 33, "closest")
  ,@("-N600j_csYpnoOceKUEq", 30, 33, "follow_stag")
  ,@("-N600kVz4mhm0haN1UoC", 34, 38, "follow_stag")
  ,@("-N600o4FxuE9h8-DIyxx", 12, 22, "random")
  ,@("-N600ut0y9Z4gOYIyiuq", 9, 17, "random")
  ,@("-N600zTc43eXAZQsIiBf", 10, 24, "random")
  ,@("-N6014T7XN-Yci_M6tLn", 3, 17, "random")
  ,@("-N601YNdKBzhyui406Tp", 6, 16, "random")
  ,@("-N601goESnQqCFcfJ-t3", 9, 22, "random")
  ,@("-N601lJmbtSwYem8DsrX", 12, 21, "random")
  ,@("-N602dhEgC76_1JFMR-Q", 6, 18, "follow_stag")
  ,@("-N605MFTS2FjoN5Wl4nd", 5, 13, "random")
  ,@("-N608zZT58boWm4nIERF", 32, 36, "follow_stag")
  ,@("-N60B-FXwfl_k6O0agPM", 9, 30, "follow_stag")
  ,@("-N60DBw3qxKIlxVPL_io", 11, 23, "random")
  ,@("-N60DUPUy4ZVeAjeD1CB", 27, 32, "follow_stag")
  ,@("-N60DjzerUzyJ9ekuEoF", 34, 25, "closest")
  ,@("-N60Duj8vhyIzhwaIZQc", 27, 29, "follow_stag")
  ,@("-N60ETyUPyoTeUbk2ErQ", 28, 23, "closest")
  ,@("-N60ErvzvYx_bVSDKugG", 9, 29, "follow_stag")
  ,@("-N60G23nZZuBppjDOofN", 31, 16, "closest")
  ,@("-N60HPdF56iDJ01ollY0", 7, 25, "follow_stag")
  ,@("-N60KUIWwsqc8zGiFzNG", 5, 18, "random")
  ,@("-N60KuKmqPd_41fLQthi", 26, 19, "closest")
  ,@("-N60Tx4zmO5_1Fgp6JlT", 27, 31, "closest")
  ,@("-N60UdZtsvjhImC1NplA", 25, 19, "closest")
  ,@("-N611yMxdy44Ur-GYHEH", 5, 17, "random")
  ,@("-N613eK641CjDIG7Sh5G", 4, 21, "random")
)

$startRow = 100
$ws.Range("A2:D2").Copy()
$ws.Range("A" + $startRow + ":D" + ($startRow + $newRows.Count - 1)).PasteSpecial(-4122)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
